# Add Job Posting with Job_Id=JD_001
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value  = "Job_Id"
$ws.Cells.Item(1,2).Value  = "Job_Title"
$ws.Cells.Item(1,3).Value  = "Job_Description"
$ws.Cells.Item(1,4).Value  = "Total_Years_Min_Exp"
$ws.Cells.Item(1,5).Value  = "Total_Years_Max_Exp"
$ws.Cells.Item(1,6).Value  = "Work_Mode"
$ws.Cells.Item(1,7).Value  = "Job_Location"
$ws.Cells.Item(1,8).Value  = "LinkedIn_Poster"
$ws.Cells.Item(1,9).Value  = "LinkedIn_Posted"
$ws.Cells.Item(1,10).Value = "Resume_received"
$ws.Cells.Item(1,11).Value = "Resume_downloaded"

# --- Data row (row 2) ---
$ws.Cells.Item(2,1).Value = "JD_001"
$ws.Cells.Item(2,2).Value = "Junior RPA Developer"
$ws.Cells.Item(2,3).Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Cells.Item(2,4).Value = 1
$ws.Cells.Item(2,5).Value = 4
$ws.Cells.Item(2,6).Value = "Remote"
$ws.Cells.Item(2,7).Value = "Hyderabad, Telangana, India"

# --- Header formatting: bold font, thin box border, centered/top-aligned ---
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
